$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-19 (Player, Position, Team)
$data = @(
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Jalen Williams", "SG,SF,PF", "Oklahoma City Thunder"),
    @("Haywood Highsmith", "SF,PF", "Miami Heat"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$wb.Save()
